$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows 2 & 3 "state" values (capitalization) ---
$ws.Range("F2").Value = "Completed"
$ws.Range("F3").Value = "Library"

# --- Widen column B ---
$ws.Columns.Item(2).ColumnWidth = 34.7

# --- Propagate row-3 formatting down to rows 4-30 so new rows match existing style ---
$ws.Range("A3:F3").Copy()
$ws.Range("A4:F30").PasteSpecial(-4122)

# --- Row 4: new "animal crossing" entry ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "animal crossing"
$ws.Range("C4").Value = "cute"
$ws.Range("D4").Value = "nintendo"
$ws.Range("E4").Value = "3ds"
$ws.Range("F4").Value = "InProgress"

# --- Rows 5-30: repeated "halo" entries ---
for ($r = 5; $r -le 30; $r++) {
    $ws.Range("B$r").Value = "halo"
    $ws.Range("C$r").Value = "nice"
    $ws.Range("D$r").Value = "someone"
    $ws.Range("E$r").Value = "xbox"
    $ws.Range("F$r").Value = "Completed"
}

# --- Match row heights for the newly added rows ---
for ($r = 4; $r -le 30; $r++) {
    $ws.Rows.Item($r).RowHeight = 16.5
}

# --- Selection moves to H10 ---
$null = $ws.Range("H10").Select()
